$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Insert a new row at row 12 ("Liked user Id" / addLike endpoint),
# which pushes the existing rows 12-20 down to 13-21.
# ------------------------------------------------------------------
$ws.Rows("12:12").Insert()

# Copy the formatting used by the other "SUB TABLE" rows (e.g. row 19,
# which is the post-insert position of the original "user Id" row) onto
# the freshly inserted row 12 so the new entry matches the workbook's
# existing look (font/fill/border per column).
$ws.Range("A19:G19").Copy()
$ws.Range("A12:G12").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows("12:12").RowHeight = 18

# Populate the new row's cells. Values are written in the same order the
# new strings were appended to the shared string table by the original
# edit ([Array], Liked user Id, PATCH, http://localhost:3005/blogs/addLike).
$ws.Range("C12").Value = "[Array]"
$ws.Range("B12").Value = "Liked user Id"
$ws.Range("G12").Value = "PATCH"
$ws.Range("F12").Value = "http://localhost:3005/blogs/addLike"
$ws.Range("D12").Value = "User Table"

# ------------------------------------------------------------------
# The hyperlinks that used to sit on F16/F17 now live on F17/F18 (one
# row down) because of the insert above; recreate them there. F11's
# hyperlink is also recreated in place so the final hyperlink order
# matches F9, F10, F17, F18, F11, F12 (new).
# ------------------------------------------------------------------
$links = @($ws.Hyperlinks)
$links[2].Delete()  # old F16 -> http://localhost:3005/comments/submitNew
$links2 = @($ws.Hyperlinks)
$links2[2].Delete()  # old F17 -> http://localhost:3005/comments/{blog Id}
$links3 = @($ws.Hyperlinks)
$links3[2].Delete()  # F11 -> http://localhost:3005/blogs/{blogId}

$ws.Hyperlinks.Add($ws.Range("F17"), "http://localhost:3005/comments/submitNew")
$ws.Hyperlinks.Add($ws.Range("F18"), "http://localhost:3005/comments/{blog Id}")
$ws.Hyperlinks.Add($ws.Range("F11"), "http://localhost:3005/blogs/{blogId}")
$ws.Hyperlinks.Add($ws.Range("F12"), "http://localhost:3005/blogs/addLike")

# Adding a hyperlink re-styles its cell with the default hyperlink theme
# (blue/underlined); restore the workbook's own "link cell" formatting
# (same visual style as F10, which was left untouched) on every cell
# whose hyperlink was just (re)created.
$ws.Range("F10").Copy()
$ws.Range("F11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F17").PasteSpecial(-4122)
$ws.Range("F18").PasteSpecial(-4122)
$ws.Range("F12").PasteSpecial(-4122)

# Match the selection left by the author on their last save.
$ws.Range("F12").Select()
